# Generate Report for Handoff
#
# Moves the localization status from "In Translation" to "Ready for
# handoff" and refreshes the handoff timestamps, across all three sheets
# (Overview, zh-cn, de-de). Also widens the "status" columns that now need
# to fit the longer "Ready for handoff" label.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ----------------------
# Overview!E2 (zh-cn status), Overview!F2 (de-de status)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# zh-cn!C2 and de-de!C2 "Status" column
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Handoff timestamps ----------------------------------------------------
# de-de handoff moved to 2016-08-12 16:47:51 (Overview!G2 mirrors de-de!H2)
$wsOverview.Range("G2").Value = "2016-08-12 16:47:51"
$wsDeDe.Range("H2").Value = "2016-08-12 16:47:51"

# zh-cn handoff moved to 2016-08-12 16:47:44
$wsZhCn.Range("H2").Value = "2016-08-12 16:47:44"

# --- Column widths: widen the "status" columns for the longer label -------
$newWidth = 16.333333333333336
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
